$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.409.28"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.845.38"
$ws.Range("E3").Value = "  +1.75%  "

$ws.Range("E4").Value = "  +1.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.98"

$ws.Range("E7").Value = "  +1.78%  "

$ws.Range("E8").Value = "  +0.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07464"
$ws.Range("E9").Value = "  +1.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8882"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.53"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.859.91"
$ws.Range("E12").Value = "  +2.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07389"
$ws.Range("E13").Value = "  +4.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.496"
$ws.Range("E14").Value = "  +2.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.30"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.592"
$ws.Range("E16").Value = "  +1.52%  "

$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008875"
$ws.Range("E18").Value = "  +2.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.014"
$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.436.15"
$ws.Range("E21").Value = "  +2.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.351"
$ws.Range("E22").Value = "  +0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.73"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.073.71"
$ws.Range("E24").Value = "  +1.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.910"
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.28"
$ws.Range("E26").Value = "  +0.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.66"
$ws.Range("E27").Value = "  +1.75%  "

$ws.Range("E28").Value = "  +0.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.292"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.06"
$ws.Range("E30").Value = "  +2.23%  "

$ws.Range("E31").Value = "  +0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7623"
$ws.Range("E32").Value = "  -0.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.180"
$ws.Range("E33").Value = "  +1.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.567"
$ws.Range("E34").Value = "  +1.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.947"
$ws.Range("E35").Value = "  +1.54%  "

$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.108"
$ws.Range("E37").Value = "  +1.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05371"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01969"
$ws.Range("E39").Value = "  +0.59%  "

$ws.Range("E40").Value = "  +2.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.324"
$ws.Range("E41").Value = "  +0.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.406"
$ws.Range("E42").Value = "  +1.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5359"
$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1669"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.560"
$ws.Range("E45").Value = "  +1.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4974"
$ws.Range("E46").Value = "  +1.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.56"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.014"
$ws.Range("E48").Value = "  +1.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.34"
$ws.Range("E49").Value = "  +1.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.685"
$ws.Range("E50").Value = "  +1.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06332"
$ws.Range("E51").Value = "  +0.90%  "
